$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "elec_config" column (C), shifting min_pk_height/min_pk_dist left.
$ws.Columns("C").Delete()

# Insert a new column for "sample_frequency" right before the (shifted) toggle_trunc column.
$ws.Columns("E").Insert()

$ws.Range("E1").Value = "sample_frequency"

$ws.Range("E2").Value = 1000
$ws.Range("E3").Value = 1000
$ws.Range("E4").Value = 1000
$ws.Range("E5").Value = 1000
$ws.Range("E6").Value = 1000
$ws.Range("E7").Value = 1000
$ws.Range("E8").Value = 1000
$ws.Range("E9").Value = 1000
$ws.Range("E10").Value = 1000
$ws.Range("E11").Value = 1000
$ws.Range("E12").Value = 1000
$ws.Range("E13").Value = 1000
$ws.Range("E14").Value = 1000

# toggle_trunc / toggle_silence become explicit boolean formulas instead of literal booleans.
$ws.Range("F2").Formula = "=FALSE()"
$ws.Range("F3").Formula = "=FALSE()"
$ws.Range("F4").Formula = "=FALSE()"
$ws.Range("F5").Formula = "=FALSE()"
$ws.Range("F6").Formula = "=FALSE()"
$ws.Range("F7").Formula = "=TRUE()"
$ws.Range("F8").Formula = "=TRUE()"
$ws.Range("F9").Formula = "=FALSE()"
$ws.Range("F10").Formula = "=FALSE()"
$ws.Range("F11").Formula = "=TRUE()"
$ws.Range("F12").Formula = "=FALSE()"
$ws.Range("F13").Formula = "=FALSE()"
$ws.Range("F14").Formula = "=FALSE()"

$ws.Range("I2").Formula = "=TRUE()"
$ws.Range("I3").Formula = "=TRUE()"
$ws.Range("I4").Formula = "=FALSE()"
$ws.Range("I5").Formula = "=FALSE()"
$ws.Range("I6").Formula = "=TRUE()"
$ws.Range("I7").Formula = "=FALSE()"
$ws.Range("I8").Formula = "=FALSE()"
$ws.Range("I9").Formula = "=TRUE()"
$ws.Range("I10").Formula = "=TRUE()"
$ws.Range("I11").Formula = "=FALSE()"
$ws.Range("I12").Formula = "=FALSE()"
$ws.Range("I13").Formula = "=FALSE()"
$ws.Range("I14").Formula = "=TRUE()"

$ws.Range("E15").Value = ""
$ws.Range("E16").Value = ""

$ws.Range("E20").Select()
